$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6 currently holds only the label "grandes regiões e unidades da federação"
# (no numeric data). Remove this entire row so that every row below it
# (7..37) shifts up by one, and the now-unused shared string is dropped.
$ws.Rows(6).Delete()
